# Update countries & provincias Spain
# Refresh the COVID-19 snapshot on the "Pais" sheet: the header timestamp
# moves from 11:22 to 11:52, several countries' case counts are updated,
# and a handful of countries change rank (their row shifts down one
# position as a new/updated country overtakes them), which also changes
# which country name sits in that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---------------------------------------------------
$ws.Cells.Item(1,1).Value = "Datos actualizados a 15 de Abril de 2020 a las 11:52"

# --- Helper: write a full data row (country label + 7 numeric columns) --
function Set-PaisRow($row, $pais, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $casosTotales
    $ws.Cells.Item($row, 3).Value = $nuevosCasos
    $ws.Cells.Item($row, 4).Value = $casosActivos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $casosCriticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# --- Rows whose numbers were refreshed (country keeps its row) ----------
Set-PaisRow 5  "España"    177633 3573 70853 88201 7371 324 18579
Set-PaisRow 16 "Suiza"     26023  87   14700 10133 386  16  1190
Set-PaisRow 31 "Rumania"   6879   0    1051  5466  241  11  362
Set-PaisRow 41 "Malasia"   5072   85   2647  2342  56   1   83
Set-PaisRow 51 "Finlandia" 3237   76   300   2873  74   0   64

# --- Kuwait overtakes Estonia / Irak / Nueva Zelanda (rows 65-68 shift) --
Set-PaisRow 65 "Kuwait"        1405 50 206 1196 31 0 3
Set-PaisRow 66 "Estonia"       1400 27 117 1248 10 4 35
Set-PaisRow 67 "Irak"          1400 0  766 556  0  0  78
Set-PaisRow 68 "Nueva Zelanda" 1386 20 728 649  3  0  9

# --- Eslovaquia overtakes Camerun (rows 80-81 shift) --------------------
Set-PaisRow 80 "Eslovaquia" 863 28 113 748 5 0 2
Set-PaisRow 81 "Camerun"    848 0  130 704 0 0 14

# --- Tanzania overtakes Liberia / Bermudas / Polinesia Francesa / Uganda
#     / Islas Caimanes (rows 146-151 shift) -------------------------------
Set-PaisRow 146 "Tanzania"            59 6 7  49 0 0 3
Set-PaisRow 147 "Liberia"             59 0 4  49 0 0 6
Set-PaisRow 148 "Bermudas"            57 0 30 22 3 0 5
Set-PaisRow 149 "Polinesia Francesa"  55 0 0  55 1 0 0
Set-PaisRow 150 "Uganda"              55 0 8  47 0 0 0
Set-PaisRow 151 "Islas Caimanes"      54 0 6  47 3 0 1
